# Insert a new data row before row 120 (pushing existing rows 120-147 down
# to 121-148) and populate it with a new weekly price observation for
# "Zapallo italiano" at Feria Lagunitas de Puerto Montt.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 120..147 down by one to make room for the new record.
$ws.Rows("120:120").Insert()

# Populate the newly inserted row 120 with the new observation.
$ws.Cells.Item(120, 1).Value  = 4
$ws.Cells.Item(120, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(120, 3).Value  = "Los Lagos"
$ws.Cells.Item(120, 4).Value  = 44511
$ws.Cells.Item(120, 5).Value  = 10
$ws.Cells.Item(120, 6).Value  = 100112032
$ws.Cells.Item(120, 7).Value  = "Zapallo italiano"
$ws.Cells.Item(120, 8).Value  = "Sin especificar"
$ws.Cells.Item(120, 9).Value  = "Primera"
$ws.Cells.Item(120, 10).Value = 120
$ws.Cells.Item(120, 11).Value = 11000
$ws.Cells.Item(120, 12).Value = 12000
$ws.Cells.Item(120, 13).Value = 11333
$ws.Cells.Item(120, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(120, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(120, 16).Value = 227
$ws.Cells.Item(120, 17).Value = 50
$ws.Cells.Item(120, 18).Value = "Hortaliza"
